$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.245863437652588
$ws.Range("B1").Value = 1.796949982643127
$ws.Range("C1").Value = 4.695855140686035
$ws.Range("D1").Value = 0.7894878387451172
$ws.Range("E1").Value = 0.7940481305122375
